$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 8 (reporting period) data updates: move from Q4 2021 to Q1 2022
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 2022
$ws.Range("B8").Value = (Get-Date -Year 2022 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C8").Value = (Get-Date -Year 2022 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K8").Value = (Get-Date -Year 2022 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L8").Value = (Get-Date -Year 2022 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)

# Updated note text (M8) -- the COVID note was rewritten
$ws.Range("M8").Value = "Debido a contingencia por COVID-19, el regreso a clases ha sido escalonado, por las características de las becas Institucionales, se otorgarán cuando todos las y los estudiantes regresen a clases presenciales, existiendo en la institución 2 tipos de becas: Alimenticia y De Representación Institucional de acuerdo a las reglas de operación publicadas el el POEH el día 28-feb-22."

# J8 "area responsable" cell gets a distinct font (Calibri, explicit black) and
# loses its inherited left alignment / date number format remnants.
$jr = $ws.Range("J8")
$jr.Style = "Normal"
$jr.Font.Name = "Calibri"
$jr.Font.Size = 11
$jr.Font.Color = 0
$jr.Borders.LineStyle = 1

# Row 8 height shrank slightly with the new note text
$ws.Rows.Item(8).RowHeight = 98.25

# ---------------------------------------------------------------------------
# Row 3 formatting: headers G3 (Tipo de programa) now wraps, as do H3:I3
# ---------------------------------------------------------------------------
$ws.Range("G3").WrapText = $true
$ws.Range("H3:I3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 39.75

# ---------------------------------------------------------------------------
# Column M got wider to fit the longer note
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 60.75

# ---------------------------------------------------------------------------
# Data validation ranges shrank from row 72 to row 57
# ---------------------------------------------------------------------------
$ws.Range("D8:D72").Validation.Delete()
$ws.Range("E8:E72").Validation.Delete()

$ws.Range("D8:D57").Validation.Add(3, 1, 1, "=Hidden_13")
$ws.Range("D8:D57").Validation.ShowInput = $false
$ws.Range("D8:D57").Validation.ShowError = $true

$ws.Range("E8:E57").Validation.Add(3, 1, 1, "=Hidden_24")
$ws.Range("E8:E57").Validation.ShowInput = $false
$ws.Range("E8:E57").Validation.ShowError = $true

# ---------------------------------------------------------------------------
# View settings: zoom 95 -> 90, selection moved up one row, print paper size
# reset to the printer default
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("M11").Select()
$ws.PageSetup.PaperSize = $null
